$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program_choosing")

$ws.Range("A5").Value = "Uni_Mannheim_MGM"
$ws.Range("B5").Value = "Yes"

$ws.Range("A6").Value = "Uni_Magdeburg_Finalcial_Economics"
$ws.Range("B6").Value = "Yes"

$ws.Range("A7").Select()
